$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "Batch" sheet (it is no longer the active tab) ---
$ws2 = $wb.Worksheets.Item("Batch")
[void]$ws2.Range("A2").Select()

# --- Add the new "Program" worksheet at the end of the workbook ---
$ws1 = $wb.Worksheets.Item("Login")
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Program"

# Copy the header style (border + left-aligned black font) used by column A headers
# on the "Login" sheet so the new sheet reuses the existing style record.
[void]$ws1.Range("A1").Copy()
[void]$newSheet.Range("A1").PasteSpecial(-4122)

# --- Populate header row ---
$newSheet.Range("A1").Value = "DataInput"
$newSheet.Range("B1").Value = "ProgramName"
$newSheet.Range("C1").Value = "ProgramDescription"
$newSheet.Range("D1").Value = "SearchCreatedProgramName"

# --- Populate data rows ---
$newSheet.Range("A2").Value = "validateTextbox"
$newSheet.Range("B2").Value = "Playwright"
$newSheet.Range("C2").Value = "Javascript"

$newSheet.Range("A3").Value = "validInput"
$newSheet.Range("B3").Value = "physeleniumplaywrightjavas"
$newSheet.Range("C3").Value = "playwright with javascript"

$newSheet.Range("A4").Value = "editProgram"
$newSheet.Range("B4").Value = "PlaywrightJavascriptJavaa"
$newSheet.Range("C4").Value = "playwright with javascript updating now"
$newSheet.Range("D4").Value = "PlaywrightJavascriptAxcf"

$newSheet.Range("A5").Value = "deleteProgram"

# --- Set column widths to match the authored sheet ---
$newSheet.Range("A1").ColumnWidth = 14.6328125
$newSheet.Range("B1").ColumnWidth = 26.08984375
$newSheet.Range("C1").ColumnWidth = 27.1796875
$newSheet.Range("D1").ColumnWidth = 25.26953125

# --- Make the new sheet the active tab, with C4 selected ---
[void]$newSheet.Range("C4").Select()
